$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 05:12"

# Apply country-table updates (reordered countries due to new case totals)
$ws.Range("D44").Value = 50
$ws.Range("E44").Value = 657
$ws.Range("E50").Value = 573
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 8
$ws.Range("A104").Value = "Camerun"
$ws.Range("B104").Value = 88
$ws.Range("C104").Value = 13
$ws.Range("D104").Value = 2
$ws.Range("E104").Value = 85
$ws.Range("A105").Value = "Estado de Palestina"
$ws.Range("D105").Value = 17
$ws.Range("E105").Value = 68
$ws.Range("F105").Value = 0
$ws.Range("H105").Value = 1
$ws.Range("A106").Value = "Bielorrusia"
$ws.Range("B106").Value = 86
$ws.Range("D106").Value = 29
$ws.Range("E106").Value = 57
$ws.Range("F106").Value = 2
$ws.Range("H106").Value = 0
$ws.Range("A107").Value = "Martinica"
$ws.Range("E107").Value = 80
$ws.Range("F107").Value = 12
$ws.Range("H107").Value = 1
$ws.Range("A108").Value = "Mauricio"
$ws.Range("B108").Value = 81
$ws.Range("D108").Value = 0
$ws.Range("E108").Value = 79
$ws.Range("H108").Value = 2
$ws.Range("A109").Value = "Georgia"
$ws.Range("B109").Value = 79
$ws.Range("D109").Value = 11
$ws.Range("E109").Value = 68
$ws.Range("F109").Value = 1
$ws.Range("A110").Value = "Uzbekistan"
$ws.Range("D110").Value = 0
$ws.Range("E110").Value = 75
$ws.Range("F110").Value = 4
$ws.Range("H110").Value = 0
$ws.Range("A143").Value = "Nueva Caledonia"
$ws.Range("A144").Value = "Uganda"
$ws.Range("A150").Value = "Mongolia"
$ws.Range("A151").Value = "San Martin (Parte Francesa)"
$ws.Range("A152").Value = "Republica de Yibuti"
$ws.Range("A156").Value = "Surinam"
$ws.Range("A157").Value = "Haiti"
$ws.Range("A160").Value = "Mozambique"
$ws.Range("A161").Value = "Granada"
$ws.Range("A162").Value = "Seychelles"
$ws.Range("A163").Value = "Antigua y Barbuda"
$ws.Range("A166").Value = "Laos"
$ws.Range("A169").Value = "Eritrea"
$ws.Range("A171").Value = "Fiyi"
$ws.Range("A172").Value = "Birmania"
$ws.Range("A173").Value = "Montserrat"
$ws.Range("A174").Value = "Siria"
$ws.Range("A175").Value = "Guyana"
$ws.Range("A176").Value = "Cabo Verde"
$ws.Range("A183").Value = "Republica del Chad"
$ws.Range("A184").Value = "San Bartolome"
$ws.Range("A186").Value = "Liberia"
$ws.Range("A187").Value = "San Martin (Parte Holandesa)"
$ws.Range("A189").Value = "Gambia"
$ws.Range("D189").Value = 0
$ws.Range("H189").Value = 1
$ws.Range("A190").Value = "Santa Lucia"
$ws.Range("D190").Value = 1
$ws.Range("H190").Value = 0
